# Hjemme passive tweaks lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header counts changed (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON"): B2, D2, E2 values removed; C2 value updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 6.0407726312651349
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 ("STR"): B3:E3 values updated
$ws.Range("B3").Value = 6.3751365426387139
$ws.Range("C3").Value = 7.7340340866256723
$ws.Range("D3").Value = 9.9670937305617571
$ws.Range("E3").Value = 4.4125299748362021

# Selection narrowed from B1:AY3 to B1:E3
$ws.Range("B1:E3").Select()
